$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy style from existing header cell (A1) to new header cells
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122) # xlPasteFormats

# Boolean values for F2:H21
$values = @(
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,1,0),
    @(0,0,0)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 6).Value = [bool]($values[$i][0])
    $ws.Cells.Item($row, 7).Value = [bool]($values[$i][1])
    $ws.Cells.Item($row, 8).Value = [bool]($values[$i][2])
}
